$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header): extend with P1=14, Q1=15, copying style (bold/border/center) from O1 ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Rows 2-25: add new columns P and Q, both filled with 2 (no special style) ---
$ws.Range("P2:Q25").Value = 2

# --- Rows 2-25: flip values in columns I, K, M, O (1<->2) ---
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I = 9
    $ws.Cells.Item($r, 9).Value = 3 - $iVal

    $kVal = $ws.Cells.Item($r, 11).Value2  # column K = 11
    $ws.Cells.Item($r, 11).Value = 3 - $kVal

    $mVal = $ws.Cells.Item($r, 13).Value2  # column M = 13
    $ws.Cells.Item($r, 13).Value = 3 - $mVal

    $oVal = $ws.Cells.Item($r, 15).Value2  # column O = 15
    $ws.Cells.Item($r, 15).Value = 3 - $oVal
}

# clear clipboard marching ants (best effort)
$excel.CutCopyMode = $false
